$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 405365.12
$ws.Cells.Item(17, 10).Value = 420588.66
$ws.Cells.Item(17, 12).Value = 1261765.98
$ws.Cells.Item(17, 14).Value = -1262101.98
$ws.Cells.Item(32, 8).Value = 4264.923
$ws.Cells.Item(32, 9).Value = 3000
$ws.Cells.Item(32, 10).Value = 4370.3335
$ws.Cells.Item(32, 11).Value = 3000
$ws.Cells.Item(32, 12).Value = 4370.3335
$ws.Cells.Item(32, 13).Value = -2674
$ws.Cells.Item(32, 14).Value = -5022.3335
$ws.Cells.Item(33, 8).Value = 68576
$ws.Cells.Item(33, 9).Value = 107419.29
$ws.Cells.Item(33, 10).Value = 600.25
$ws.Cells.Item(33, 11).Value = 107419.29
$ws.Cells.Item(33, 12).Value = 600.25
$ws.Cells.Item(33, 13).Value = -107190.29
$ws.Cells.Item(33, 14).Value = -1058.25
$ws.Cells.Item(40, 8).Value = 2672.1667
$ws.Cells.Item(40, 9).Value = 2487.5
$ws.Cells.Item(40, 10).Value = 2819.9
$ws.Cells.Item(40, 11).Value = 2487.5
$ws.Cells.Item(40, 12).Value = 2819.9
$ws.Cells.Item(40, 13).Value = -2312.5
$ws.Cells.Item(40, 14).Value = -3169.9
$ws.Cells.Item(51, 8).Value = 5308.3
$ws.Cells.Item(51, 9).Value = 8000
$ws.Cells.Item(51, 10).Value = 5166.6313
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 5166.6313
$ws.Cells.Item(51, 13).Value = -7516
$ws.Cells.Item(51, 14).Value = -6134.6313
$ws.Cells.Item(76, 8).Value = 21458334
$ws.Cells.Item(76, 9).Value = 499999
$ws.Cells.Item(76, 10).Value = 28444444
$ws.Cells.Item(76, 11).Value = 499999
$ws.Cells.Item(76, 12).Value = 28444444
$ws.Cells.Item(76, 13).Value = -499684
$ws.Cells.Item(76, 14).Value = -28445074
$ws.Cells.Item(79, 8).Value = 21458334
$ws.Cells.Item(79, 9).Value = 499999
$ws.Cells.Item(79, 10).Value = 28444444
$ws.Cells.Item(79, 11).Value = 499999
$ws.Cells.Item(79, 12).Value = 28444444
$ws.Cells.Item(79, 13).Value = -498907
$ws.Cells.Item(79, 14).Value = -28446628
$ws.Cells.Item(86, 8).Value = 22248886
$ws.Cells.Item(86, 9).Value = 9991.333000000001
$ws.Cells.Item(86, 10).Value = 33368334
$ws.Cells.Item(86, 11).Value = 9991.333000000001
$ws.Cells.Item(86, 12).Value = 33368334
$ws.Cells.Item(86, 13).Value = -8868.333000000001
$ws.Cells.Item(86, 14).Value = -33370580
$ws.Cells.Item(89, 8).Value = 22248886
$ws.Cells.Item(89, 9).Value = 9991.333000000001
$ws.Cells.Item(89, 10).Value = 33368334
$ws.Cells.Item(89, 11).Value = 49956.665
$ws.Cells.Item(89, 12).Value = 166841670
$ws.Cells.Item(89, 13).Value = -44340.665
$ws.Cells.Item(89, 14).Value = -166852902
$ws.Cells.Item(129, 8).Value = 1197
$ws.Cells.Item(129, 9).Value = 942.3333
$ws.Cells.Item(129, 10).Value = 1579
$ws.Cells.Item(129, 11).Value = 2826.9999
$ws.Cells.Item(129, 12).Value = 4737
$ws.Cells.Item(129, 13).Value = 2173.0001
$ws.Cells.Item(129, 14).Value = -14737
$ws.Cells.Item(131, 8).Value = 20533.166
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 4424.1333
$ws.Cells.Item(132, 9).Value = 4110.24
$ws.Cells.Item(132, 11).Value = 12330.72
$ws.Cells.Item(132, 13).Value = -9800.719999999999
$ws.Cells.Item(141, 8).Value = 5056.2856
$ws.Cells.Item(141, 9).Value = 5119
$ws.Cells.Item(141, 10).Value = 4899.5
$ws.Cells.Item(141, 11).Value = 15357
$ws.Cells.Item(141, 12).Value = 14698.5
$ws.Cells.Item(141, 13).Value = -10177
$ws.Cells.Item(141, 14).Value = -25058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3472.1738
$ws.Cells.Item(32, 9).Value = 3472.1738
$ws.Cells.Item(32, 11).Value = 3472.1738
$ws.Cells.Item(32, 13).Value = -3185.1738
$ws.Cells.Item(110, 8).Value = 64917.562
$ws.Cells.Item(110, 9).Value = 44056.832
$ws.Cells.Item(110, 10).Value = 127499.75
$ws.Cells.Item(110, 11).Value = 44056.832
$ws.Cells.Item(110, 12).Value = 127499.75
$ws.Cells.Item(110, 13).Value = -42011.832
$ws.Cells.Item(110, 14).Value = -131589.75
$ws.Cells.Item(122, 8).Value = 18521858
$ws.Cells.Item(122, 9).Value = 27780286
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 83340858
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -83338408
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(132, 8).Value = 50002684
$ws.Cells.Item(132, 9).Value = 71430890
$ws.Cells.Item(132, 10).Value = 3535.5
$ws.Cells.Item(132, 11).Value = 214292670
$ws.Cells.Item(132, 12).Value = 10606.5
$ws.Cells.Item(132, 13).Value = -214290140
$ws.Cells.Item(132, 14).Value = -15666.5
$ws.Cells.Item(135, 8).Value = 31101
$ws.Cells.Item(135, 10).Value = 31101
$ws.Cells.Item(135, 12).Value = 31101
$ws.Cells.Item(135, 14).Value = -41241
$ws.Cells.Item(137, 8).Value = 59896
$ws.Cells.Item(137, 10).Value = 59896
$ws.Cells.Item(137, 12).Value = 59896
$ws.Cells.Item(137, 14).Value = -70096

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 26318968
$ws.Cells.Item(86, 9).Value = 35717670
$ws.Cells.Item(86, 11).Value = 35717670
$ws.Cells.Item(86, 13).Value = -35716547
$ws.Cells.Item(89, 8).Value = 26318968
$ws.Cells.Item(89, 9).Value = 35717670
$ws.Cells.Item(89, 11).Value = 178588350
$ws.Cells.Item(89, 13).Value = -178582734
$ws.Cells.Item(105, 8).Value = 1835.24
$ws.Cells.Item(105, 9).Value = 1758.1111
$ws.Cells.Item(105, 11).Value = 1758.1111
$ws.Cells.Item(105, 13).Value = -11.11110000000008
$ws.Cells.Item(107, 8).Value = 26878.7
$ws.Cells.Item(107, 9).Value = 22296.928
$ws.Cells.Item(107, 11).Value = 22296.928
$ws.Cells.Item(107, 13).Value = -20376.928

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(75, 8).Value = 12222
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 12222
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 4359.625
$ws.Cells.Item(86, 10).Value = 5526.3335
$ws.Cells.Item(86, 12).Value = 5526.3335
$ws.Cells.Item(86, 14).Value = -7772.3335
$ws.Cells.Item(89, 8).Value = 4359.625
$ws.Cells.Item(89, 10).Value = 5526.3335
$ws.Cells.Item(89, 12).Value = 27631.6675
$ws.Cells.Item(89, 14).Value = -38863.6675
$ws.Cells.Item(107, 8).Value = 2024.9584
$ws.Cells.Item(107, 9).Value = 1901.9286
$ws.Cells.Item(107, 10).Value = 2197.2
$ws.Cells.Item(107, 11).Value = 1901.9286
$ws.Cells.Item(107, 12).Value = 2197.2
$ws.Cells.Item(107, 13).Value = 18.07140000000004
$ws.Cells.Item(107, 14).Value = -6037.2
$ws.Cells.Item(122, 8).Value = 2647.9375
$ws.Cells.Item(122, 9).Value = 2262.75
$ws.Cells.Item(122, 10).Value = 3033.125
$ws.Cells.Item(122, 11).Value = 6788.25
$ws.Cells.Item(122, 12).Value = 9099.375
$ws.Cells.Item(122, 13).Value = -4338.25
$ws.Cells.Item(122, 14).Value = -13999.375
$ws.Cells.Item(134, 8).Value = 3278.6667
$ws.Cells.Item(134, 9).Value = 2891.6667
$ws.Cells.Item(134, 11).Value = 8675.000100000001
$ws.Cells.Item(134, 13).Value = -6140.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 802.3
$ws.Cells.Item(34, 10).Value = 1175
$ws.Cells.Item(34, 12).Value = 3525
$ws.Cells.Item(34, 14).Value = -3693
$ws.Cells.Item(98, 8).Value = 7595.125
$ws.Cells.Item(98, 9).Value = 14328.5
$ws.Cells.Item(98, 10).Value = 861.75
$ws.Cells.Item(98, 11).Value = 42985.5
$ws.Cells.Item(98, 12).Value = 2585.25
$ws.Cells.Item(98, 13).Value = -41487.5
$ws.Cells.Item(98, 14).Value = -5581.25
$ws.Cells.Item(110, 8).Value = 15000
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 13).ClearContents()
$ws.Cells.Item(112, 8).Value = 1799.6666
$ws.Cells.Item(112, 9).Value = 1799.6666
$ws.Cells.Item(112, 11).Value = 5398.9998
$ws.Cells.Item(112, 13).Value = -4290.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6021.0415
$ws.Cells.Item(70, 9).Value = 4248
$ws.Cells.Item(70, 11).Value = 4248
$ws.Cells.Item(70, 13).Value = -3978
$ws.Cells.Item(73, 8).Value = 6021.0415
$ws.Cells.Item(73, 9).Value = 4248
$ws.Cells.Item(73, 11).Value = 4248
$ws.Cells.Item(73, 13).Value = -3312
$ws.Cells.Item(113, 8).Value = 2226.5454
$ws.Cells.Item(113, 9).Value = 1938.375
$ws.Cells.Item(113, 10).Value = 2995
$ws.Cells.Item(113, 11).Value = 1938.375
$ws.Cells.Item(113, 12).Value = 2995
$ws.Cells.Item(113, 13).Value = 231.625
$ws.Cells.Item(113, 14).Value = -7335
$ws.Cells.Item(122, 8).Value = 2178.9688
$ws.Cells.Item(122, 9).Value = 2051.682
$ws.Cells.Item(122, 10).Value = 2459
$ws.Cells.Item(122, 11).Value = 6155.045999999999
$ws.Cells.Item(122, 12).Value = 7377
$ws.Cells.Item(122, 13).Value = -3705.045999999999
$ws.Cells.Item(122, 14).Value = -12277
$ws.Cells.Item(132, 8).Value = 3690.5676
$ws.Cells.Item(132, 9).Value = 2867.6924
$ws.Cells.Item(132, 11).Value = 8603.0772
$ws.Cells.Item(132, 13).Value = -6073.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2548.2856
$ws.Cells.Item(46, 9).Value = 968.5
$ws.Cells.Item(46, 10).Value = 3180.2
$ws.Cells.Item(46, 11).Value = 968.5
$ws.Cells.Item(46, 12).Value = 3180.2
$ws.Cells.Item(46, 13).Value = -780.5
$ws.Cells.Item(46, 14).Value = -3556.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 1499
$ws.Cells.Item(12, 9).Value = 1499
$ws.Cells.Item(12, 11).Value = 1499
$ws.Cells.Item(12, 13).Value = -1357
$ws.Cells.Item(23, 8).Value = 1608.5454
$ws.Cells.Item(23, 9).Value = 1689.4
$ws.Cells.Item(23, 10).Value = 800
$ws.Cells.Item(23, 11).Value = 1689.4
$ws.Cells.Item(23, 12).Value = 800
$ws.Cells.Item(23, 13).Value = -1460.4
$ws.Cells.Item(23, 14).Value = -1258
$ws.Cells.Item(81, 8).Value = 7412773.5
$ws.Cells.Item(81, 10).Value = 22233550
$ws.Cells.Item(81, 12).Value = 44467100
$ws.Cells.Item(81, 14).Value = -44469222
$ws.Cells.Item(84, 8).Value = 7412773.5
$ws.Cells.Item(84, 10).Value = 22233550
$ws.Cells.Item(84, 12).Value = 222335500
$ws.Cells.Item(84, 14).Value = -222346108
$ws.Cells.Item(96, 8).Value = 2771.4285
$ws.Cells.Item(96, 9).Value = 2250
$ws.Cells.Item(96, 10).Value = 2980
$ws.Cells.Item(96, 11).Value = 2250
$ws.Cells.Item(96, 12).Value = 2980
$ws.Cells.Item(96, 13).Value = -877
$ws.Cells.Item(96, 14).Value = -5726
$ws.Cells.Item(106, 8).Value = 25000
$ws.Cells.Item(106, 10).Value = 25000
$ws.Cells.Item(106, 12).Value = 25000
$ws.Cells.Item(106, 14).Value = -27524
$ws.Cells.Item(122, 8).Value = 2191.8333
$ws.Cells.Item(122, 9).Value = 2290
$ws.Cells.Item(122, 10).Value = 1818.8
$ws.Cells.Item(122, 11).Value = 6870
$ws.Cells.Item(122, 12).Value = 5456.4
$ws.Cells.Item(122, 13).Value = -4420
$ws.Cells.Item(122, 14).Value = -10356.4
$ws.Cells.Item(126, 8).Value = 2227.8572
$ws.Cells.Item(126, 10).Value = 2531.6667
$ws.Cells.Item(126, 12).Value = 7595.000100000001
$ws.Cells.Item(126, 14).Value = -12535.0001
